$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 105, shifting existing rows 105..126 down to 106..127
$ws.Rows.Item(105).Insert()

# Populate the new row 105 with the new price-report record
$ws.Range("A105").Value = 10
$ws.Range("B105").Value = "Vega Modelo de Temuco"
$ws.Range("C105").Value = "La Araucanía"
$ws.Range("D105").Value = 44722
$ws.Range("E105").Value = 9
$ws.Range("F105").Value = "Fruta"
$ws.Range("G105").Value = 100104
$ws.Range("H105").Value = "Frutos de pepita"
$ws.Range("I105").Value = 100104001
$ws.Range("J105").Value = "Granada"
$ws.Range("K105").Value = "Wonderfull"
$ws.Range("L105").Value = "Primera"
$ws.Range("M105").Value = 65
$ws.Range("N105").Value = 15000
$ws.Range("O105").Value = 15000
$ws.Range("P105").Value = 15000
$ws.Range("Q105").Value = "$/bandeja 15 kilos granel"
$ws.Range("R105").Value = "Provincia de Limarí"
$ws.Range("S105").Value = 1000
$ws.Range("T105").Value = 15

# Match the D105 cell's number/date style to the rest of the date column
$ws.Range("D105").NumberFormat = $ws.Range("D106").NumberFormat
